# "Created Program Feature file Addon"
# Adds two new menu/header entries ("Manage Class" and "Dashboard") to the
# Login sheet's H column (rows 5 and 6), and selects cell H2 on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# New data values for the table (these also register two new shared strings:
# "Manage Class" and "Dashboard").
$ws.Range("H5").Value = "Manage Class"
$ws.Range("H6").Value = "Dashboard"

# Make sure the Login sheet is active and select H2, matching the saved
# selection state recorded in the sheet view.
$ws.Activate()
$ws.Range("H2").Select()
